# Update the "Metadata" sheet of the CodeSystem workbook:
#  - Status goes from "draft" to "active"
#  - Date is bumped to the new publish timestamp
#  - Case Sensitive is now declared ("true")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B6").Value = "active"
$ws.Range("B8").Value = "2024-12-16T14:50:05-03:00"

# Case Sensitive is now declared as "true". Writing the literal word "true"
# directly would be auto-coerced to a Boolean, so compute it as text via a
# formula and then collapse the formula down to its plain text result.
$caseSensitive = $ws.Range("B17")
$caseSensitive.Formula = "=T(""true"")"
$caseSensitive.Copy()
$caseSensitive.PasteSpecial(-4163)  # xlPasteValues
